# Applies the "Corrected excel sheets for application fix issues" edit:
#  - Summary sheet: a few corrected totals
#  - Repayment schedule sheet: schedule recalculated (dates/days/interest/
#    balances shifted/updated), the now-unused "O" column cleared out on the
#    data rows, and the final row's paid/outstanding totals bumped to 1075
#  - Transactions sheet: corrected transaction IDs
# Selections on each sheet are also moved to match where the author last
# clicked while reviewing the fix.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Activate()

$summary.Range("F2").Value = 0
$summary.Range("A3").Value = 720.4
$summary.Range("E3").Value = 520.4

$summary.Range("C4").Select()

# ---------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------
$sched = $wb.Worksheets.Item("Repayment schedule")
$sched.Activate()

# The "O" column (duplicate "Over Due" figures) is no longer populated on
# the data rows - fully clear those cells (style included) instead of just
# blanking the value.
$sched.Range("O3:O13").Clear()
# Likewise the stray P2 cell (no value, just inherited style) is removed.
$sched.Range("P2").Clear()

# Row 4
$sched.Range("B4").Value = 31
$sched.Range("C4").Value = 42095
$sched.Range("F4").Value = 872.19
$sched.Range("G4").Value = 8363.27
$sched.Range("H4").Value = 92.35

# Row 5
$sched.Range("B5").Value = 30
$sched.Range("C5").Value = 42125
$sched.Range("F5").Value = 880.91
$sched.Range("G5").Value = 7482.36
$sched.Range("H5").Value = 83.63

# Row 6
$sched.Range("B6").Value = 31
$sched.Range("C6").Value = 42156
$sched.Range("F6").Value = 889.72
$sched.Range("G6").Value = 6592.64
$sched.Range("H6").Value = 74.819999999999993

# Row 7
$sched.Range("B7").Value = 30
$sched.Range("C7").Value = 42186
$sched.Range("F7").Value = 898.61
$sched.Range("G7").Value = 5694.03
$sched.Range("H7").Value = 65.930000000000007

# Row 8
$sched.Range("B8").Value = 31
$sched.Range("C8").Value = 42217
$sched.Range("F8").Value = 907.6
$sched.Range("G8").Value = 4786.43
$sched.Range("H8").Value = 56.94

# Row 9
$sched.Range("C9").Value = 42248
$sched.Range("F9").Value = 916.68
$sched.Range("G9").Value = 3869.75
$sched.Range("H9").Value = 47.86

# Row 10
$sched.Range("B10").Value = 30
$sched.Range("C10").Value = 42278
$sched.Range("F10").Value = 925.84
$sched.Range("G10").Value = 2943.91
$sched.Range("H10").Value = 38.700000000000003

# Row 11
$sched.Range("B11").Value = 31
$sched.Range("C11").Value = 42309
$sched.Range("F11").Value = 935.1
$sched.Range("G11").Value = 2008.81
$sched.Range("H11").Value = 29.44

# Row 12 - G12 switches from the plain general style to a #,##0.00 number
# format (matching the rest of the "Balance of Loan" column).
$sched.Range("B12").Value = 30
$sched.Range("C12").Value = 42339
$sched.Range("F12").Value = 944.45
$sched.Range("G12").Value = 1064.3599999999999
$sched.Range("G12").NumberFormat = "#,##0.00"
$sched.Range("H12").Value = 20.09

# Row 13 - final row: F13 also becomes #,##0.00, and the paid/outstanding
# totals (K13/P13) switch to a #,##0 number format.
$sched.Range("B13").Value = 31
$sched.Range("C13").Value = 42370
$sched.Range("F13").Value = 1064.3599999999999
$sched.Range("F13").NumberFormat = "#,##0.00"
$sched.Range("H13").Value = 10.64
$sched.Range("K13").Value = 1075
$sched.Range("K13").NumberFormat = "#,##0"
$sched.Range("P13").Value = 1075
$sched.Range("P13").NumberFormat = "#,##0"

$sched.Range("F4:F13").Select()

# ---------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------
$txns = $wb.Worksheets.Item("Transactions")
$txns.Activate()

$txns.Range("A2").Value = 6352
$txns.Range("A3").Value = 6350

$txns.Range("D3").Select()
